$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

Set-TextValue "D2" "29.861.66"
Set-TextValue "E2" "  +2.93%  "
Set-TextValue "D3" "1.868.25"
Set-TextValue "E3" "  +2.10%  "
Set-TextValue "D4" "0.9996"
Set-TextValue "E4" "  -0.21%  "
Set-TextValue "D5" "246.93"
Set-TextValue "E5" "  +3.10%  "
Set-TextValue "D6" "0.7010"
Set-TextValue "E6" "  +2.43%  "
Set-TextValue "D7" "0.9998"
Set-TextValue "E7" "  -0.19%  "
Set-TextValue "D8" "0.07786"
Set-TextValue "E8" "  +1.95%  "
Set-TextValue "D9" "0.3089"
Set-TextValue "E9" "  +2.59%  "
Set-TextValue "D10" "23.94"
Set-TextValue "E10" "  +2.53%  "
Set-TextValue "D11" "0.07851"
Set-TextValue "E11" "  +1.39%  "
Set-TextValue "D12" "5.193"
Set-TextValue "E12" "  +2.79%  "
Set-TextValue "D13" "1.863.88"
Set-TextValue "E13" "  +1.40%  "
Set-TextValue "D14" "92.90"
Set-TextValue "E14" "  +2.96%  "
Set-TextValue "D15" "0.6971"
Set-TextValue "E15" "  +3.58%  "
Set-TextValue "D16" "6.642"
Set-TextValue "E16" "  +3.28%  "
Set-TextValue "D17" "29.843.39"
Set-TextValue "E17" "  +2.83%  "
Set-TextValue "D18" "0.000008422"
Set-TextValue "E18" "  +1.82%  "
Set-TextValue "B19" "BitcoinCash"
Set-TextValue "C19" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D19" "244.50"
Set-TextValue "E19" "  +0.65%  "
Set-TextValue "B20" "WrappedliquidstakedEther2.0"
Set-TextValue "C20" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D20" "2.113.90"
Set-TextValue "E20" "  +0.44%  "
Set-TextValue "E21" "  +1.50%  "
Set-TextValue "D22" "1.000"
Set-TextValue "E22" "  -0.21%  "
Set-TextValue "D23" "7.681"
Set-TextValue "E23" "  +2.89%  "
Set-TextValue "D24" "1.000"
Set-TextValue "E24" "  -0.17%  "
Set-TextValue "D25" "0.1515"
Set-TextValue "E25" "  +3.26%  "
Set-TextValue "D26" "8.994"
Set-TextValue "E26" "  +2.83%  "
Set-TextValue "D27" "160.22"
Set-TextValue "E27" "  -0.58%  "
Set-TextValue "D28" "18.45"
Set-TextValue "E28" "  +1.90%  "
Set-TextValue "D29" "1.544"
Set-TextValue "E29" "  +0.88%  "
Set-TextValue "E30" "  +2.15%  "
Set-TextValue "E31" "  +2.04%  "
Set-TextValue "E32" "  +1.19%  "
Set-TextValue "D33" "0.05114"
Set-TextValue "E33" "  -0.01%  "
Set-TextValue "E34" "  +4.23%  "
Set-TextValue "D35" "1.942"
Set-TextValue "E35" "  +7.01%  "
Set-TextValue "E36" "  +1.60%  "
Set-TextValue "D37" "2.713"
Set-TextValue "E37" "  +0.29%  "
Set-TextValue "D38" "1.332.57"
Set-TextValue "E38" "  +9.72%  "
Set-TextValue "D39" "0.01889"
Set-TextValue "E39" "  +3.10%  "
Set-TextValue "E40" "  +1.64%  "
Set-TextValue "D41" "0.9546"
Set-TextValue "E41" "  +4.30%  "
Set-TextValue "D42" "6.075"
Set-TextValue "E42" "  +11.73%  "
Set-TextValue "D43" "107.47"
Set-TextValue "E43" "  -1.29%  "
Set-TextValue "D44" "0.9999"
Set-TextValue "E44" "  -0.16%  "
Set-TextValue "D45" "9.851"
Set-TextValue "E45" "  +4.62%  "
Set-TextValue "B46" "BabyDogeCoin"
Set-TextValue "C46" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D46" "0.00000000126"
Set-TextValue "E46" "  +5.17%  "
Set-TextValue "B47" "RocketPoolETH"
Set-TextValue "C47" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D47" "2.012.50"
Set-TextValue "E47" "  +0.30%  "
Set-TextValue "B48" "Aave"
Set-TextValue "C48" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D48" "65.58"
Set-TextValue "E48" "  +3.34%  "
Set-TextValue "B49" "RenderToken"
Set-TextValue "C49" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D49" "1.802"
Set-TextValue "E49" "  +4.38%  "
Set-TextValue "B50" "Mantle"
Set-TextValue "C50" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D50" "0.5202"
Set-TextValue "E50" "  +0.66%  "
Set-TextValue "D51" "7.049"
Set-TextValue "E51" "  +2.00%  "
